$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that Word had left at the very
#    top of the document (empty centered paragraph before the title).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Tidy up the example code:
#      during (fun world -> ActivityState.isActing ^ character.GetActivityState world) ^ chain {
#    becomes
#      during (ActivityState.isActing << character.GetActivityState world) ^ chain {
# ------------------------------------------------------------------

# 2a. Drop the "fun world -> " lambda header, leaving the
#     "ActivityState.isActing" identifier (and its spell-check wrapping)
#     untouched.
$null = $d.Content.Find.Execute("fun world -> ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2b. Swap the "^" between the two identifiers for the "<<" operator.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("isActing ^ character", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $caretStart = $rng.Start + 9
    $caretEnd = $caretStart + 1
    $caretRng = $d.Range($caretStart, $caretEnd)
    $caretRng.Text = "<<"
}

# ------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark right where the cursor was left
#    after the last edit -- between "world" and the closing ") ^ chain".
# ------------------------------------------------------------------
$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute("character.GetActivityState world", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $bmRng = $d.Range($rng2.End, $rng2.End)
    $d.Bookmarks.Add("_GoBack", $bmRng)
}
